# Reposition "TextBox 6" (the Oikonomides citation textbox) on the slide.
# Target (EMU): off x=2093596 y=10034260, ext cx=3037840 cy=1315104
# PowerPoint's COM object model works in points (1 pt = 12700 EMU), so the
# EMU offsets/extents from the OOXML diff are converted to points below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 6") {
        $shape = $candidate
        break
    }
}

$shape.Left = 164.85007874015747
$shape.Top = 790.0992126984252
$shape.Width = 239.2000046
$shape.Height = 103.55149606299213
